$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift formatting one column to the left (B:F -> A:E) without disturbing the
# <cols> width definitions (which stay anchored to columns B and C:E).
$ws.Range("B1:F6").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Write the shifted literal values into A1:E6 (App Code/T001 column is gone).
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "Deployment date"
$ws.Range("C1").Value = "Type Day"
$ws.Range("D1").Value = "Deploy Program"
$ws.Range("E1").Value = "Deplot For User"

$ws.Range("A2").Value = "2019"
$ws.Range("B2").Value = "01/01/2019"
$ws.Range("C2").Value = "H"
$ws.Range("D2").Value = "a"

$ws.Range("A3").Value = "2019"
$ws.Range("B3").Value = "01/02/2019"
$ws.Range("C3").Value = "H"

$ws.Range("A4").Value = "2019"
$ws.Range("B4").Value = "22/05/2019"
$ws.Range("C4").Value = "H"

$ws.Range("A5").Value = "2019"
$ws.Range("B5").Value = "30/07/2019"
$ws.Range("C5").Value = "H"
$ws.Range("D5").Value = "s"

$ws.Range("A6").Value = "2019"
$ws.Range("B6").Value = "30/05/2019"
$ws.Range("C6").Value = "H"

# Fully clear cells that must not exist at all after the shift (value + style).
$ws.Range("E2").Clear()
$ws.Range("D3").Clear()
$ws.Range("E3").Clear()
$ws.Range("D4").Clear()
$ws.Range("E4").Clear()
$ws.Range("E5").Clear()
$ws.Range("D6").Clear()
$ws.Range("E6").Clear()

# Drop the now-duplicate rightmost column.
$ws.Columns("F").Delete()

# Match the author's final selection.
$ws.Range("C9").Select()
